$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Job": rename the table column header from "蠢事装备" to "初始装备"
# ---------------------------------------------------------------------------
$wsJob = $wb.Worksheets.Item("Job")
$wsJob.Range("H1").Value = "初始装备"

# ---------------------------------------------------------------------------
# Sheet "~说明": drop the old column E (skill description text) and add new
# element / role tag columns E..K with the relevant job tags.
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("~说明")

# Remove the old descriptive text that used to live in column E (rows 2-9).
$wsInfo.Range("E2:E9").ClearFormats()
$wsInfo.Range("E2:E9").ClearContents()

# New header row (row 1): element types
$wsInfo.Range("E1").Value = "无"
$wsInfo.Range("F1").Value = "水"
$wsInfo.Range("G1").Value = "风"
$wsInfo.Range("H1").Value = "火"
$wsInfo.Range("I1").Value = "地"
$wsInfo.Range("J1").Value = "光"
$wsInfo.Range("K1").Value = "暗"
$wsInfo.Range("E1:K1").Font.Bold = $true
$wsInfo.Range("E1:K1").Interior.ThemeColor = 5
$wsInfo.Range("E1:K1").Interior.TintAndShade = -0.59999389629810485

# Row 3 - 战士 (Warrior)
$wsInfo.Range("G3").Value = "冲锋战"
$wsInfo.Range("H3").Value = "攻击战"

# Row 4 - 护卫 (Guard)
$wsInfo.Range("I4").Value = "防御"
$wsInfo.Range("J4").Value = "回复盾"

# Row 5 - 盗贼 (Rogue)
$wsInfo.Range("F5").Value = "下毒"
$wsInfo.Range("K5").Value = "刺杀"

# Row 6 - 射手 (Archer)
$wsInfo.Range("F6").Value = "远程强化"
$wsInfo.Range("H6").Value = "伤害强化"

# Row 7 - 游侠 (Ranger)
$wsInfo.Range("E7").Value = "武器强化"
$wsInfo.Range("G7").Value = "输出"

# Row 8 - 法师 (Mage)
$wsInfo.Range("F8").Value = "冰法"
$wsInfo.Range("H8").Value = "火法"

# Row 9 - 贤者 (Sage)
$wsInfo.Range("J9").Value = "奶"
$wsInfo.Range("K9").Value = "黑暗贤者"

# Row 10 - 诗人 (Bard)
$wsInfo.Range("E10").Value = "光环怪"
$wsInfo.Range("I10").Value = "祝福"

$wsInfo.Range("E1:K1").ColumnWidth = 7

# ---------------------------------------------------------------------------
# Update the active selections to match the saved workbook state
# ---------------------------------------------------------------------------
$wsJob.Activate()
$wsJob.Range("H2").Select()

$wsInfo.Activate()
$wsInfo.Range("H6").Select()

$wsJob.Activate()
